$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D209").Value = 44637
$ws.Range("I209").Value = "Segunda"
$ws.Range("K209").Value = 1000
$ws.Range("L209").Value = 1000
$ws.Range("M209").Value = 1000
$ws.Range("P209").Value = 1000

$ws.Range("I210").Value = "Primera"
$ws.Range("J210").Value = 2000
$ws.Range("K210").Value = 500
$ws.Range("L210").Value = 500
$ws.Range("M210").Value = 500
$ws.Range("P210").Value = 500

$ws.Range("D211").Value = 44445
$ws.Range("I211").Value = "Segunda"
$ws.Range("J211").Value = 3000
$ws.Range("K211").Value = 300
$ws.Range("L211").Value = 300
$ws.Range("M211").Value = 300
$ws.Range("P211").Value = 300

$ws.Range("D212").Value = 44524
$ws.Range("J212").Value = 4000
$ws.Range("K212").Value = 700
$ws.Range("L212").Value = 700
$ws.Range("M212").Value = 700
$ws.Range("P212").Value = 700

$ws.Range("D213").Value = 44355
$ws.Range("J213").Value = 5000
$ws.Range("K213").Value = 450
$ws.Range("L213").Value = 450
$ws.Range("M213").Value = 450
$ws.Range("P213").Value = 450

$ws.Range("D214").Value = 44294
$ws.Range("J214").Value = 4000
$ws.Range("K214").Value = 800
$ws.Range("L214").Value = 800
$ws.Range("M214").Value = 800
$ws.Range("P214").Value = 800

$ws.Range("D215").Value = 44557
$ws.Range("J215").Value = 5000
$ws.Range("K215").Value = 600
$ws.Range("L215").Value = 600
$ws.Range("M215").Value = 600
$ws.Range("P215").Value = 600

$ws.Range("D216").Value = 44264
$ws.Range("J216").Value = 3000
$ws.Range("K216").Value = 1000
$ws.Range("L216").Value = 1000
$ws.Range("M216").Value = 1000
$ws.Range("P216").Value = 1000

$ws.Range("D217").Value = 44396
$ws.Range("J217").Value = 5000
$ws.Range("K217").Value = 350
$ws.Range("L217").Value = 350
$ws.Range("M217").Value = 350
$ws.Range("P217").Value = 350

$ws.Range("D218").Value = 44279
$ws.Range("K218").Value = 800
$ws.Range("L218").Value = 800
$ws.Range("M218").Value = 800
$ws.Range("P218").Value = 800

$ws.Range("D219").Value = 44330
$ws.Range("J219").Value = 3000
$ws.Range("K219").Value = 500
$ws.Range("L219").Value = 500
$ws.Range("M219").Value = 500
$ws.Range("O219").Value = "Región del Maule"
$ws.Range("P219").Value = 500

$ws.Range("D220").Value = 44504
$ws.Range("J220").Value = 6000
$ws.Range("K220").Value = 600
$ws.Range("L220").Value = 600
$ws.Range("M220").Value = 600
$ws.Range("O220").Value = "Provincia del Elquí"
$ws.Range("P220").Value = 600

$ws.Range("D221").Value = 44301
$ws.Range("J221").Value = 3000
$ws.Range("K221").Value = 700
$ws.Range("L221").Value = 700
$ws.Range("M221").Value = 700
$ws.Range("P221").Value = 700

$ws.Range("D222").Value = 44370
$ws.Range("J222").Value = 5000
$ws.Range("K222").Value = 400
$ws.Range("L222").Value = 400
$ws.Range("M222").Value = 400
$ws.Range("P222").Value = 400

$ws.Range("I223").Value = "Primera"
$ws.Range("K223").Value = 500
$ws.Range("L223").Value = 500
$ws.Range("M223").Value = 500
$ws.Range("P223").Value = 500

$ws.Range("D224").Value = 44413
$ws.Range("I224").Value = "Segunda"
$ws.Range("J224").Value = 3000
$ws.Range("K224").Value = 350
$ws.Range("L224").Value = 350
$ws.Range("M224").Value = 350
$ws.Range("P224").Value = 350

$ws.Range("D225").Value = 44272
$ws.Range("J225").Value = 2000
$ws.Range("K225").Value = 800
$ws.Range("L225").Value = 800
$ws.Range("M225").Value = 800
$ws.Range("P225").Value = 800

$ws.Range("D226").Value = 44214
$ws.Range("J226").Value = 3000
$ws.Range("K226").Value = 900
$ws.Range("L226").Value = 900
$ws.Range("M226").Value = 900
$ws.Range("P226").Value = 900

$ws.Range("D227").Value = 44312

$ws.Range("D228").Value = 44399
$ws.Range("K228").Value = 400
$ws.Range("L228").Value = 400
$ws.Range("M228").Value = 400
$ws.Range("P228").Value = 400

$ws.Range("D229").Value = 44543
$ws.Range("J229").Value = 5000
$ws.Range("K229").Value = 500
$ws.Range("L229").Value = 500
$ws.Range("M229").Value = 500
$ws.Range("P229").Value = 500

$ws.Range("D230").Value = 44167
$ws.Range("K230").Value = 700
$ws.Range("L230").Value = 700
$ws.Range("M230").Value = 700
$ws.Range("P230").Value = 700

$ws.Range("D231").Value = 44277
$ws.Range("H231").Value = "Crespo record"
$ws.Range("J231").Value = 3000

$ws.Range("D232").Value = 44258
$ws.Range("H232").Value = "Copenhague"
$ws.Range("J232").Value = 2000
$ws.Range("K232").Value = 800
$ws.Range("L232").Value = 800
$ws.Range("M232").Value = 800
$ws.Range("P232").Value = 800

$ws.Range("D233").Value = 44390
$ws.Range("K233").Value = 400
$ws.Range("L233").Value = 400
$ws.Range("M233").Value = 400
$ws.Range("P233").Value = 400

$ws.Range("D234").Value = 44349
$ws.Range("J234").Value = 5000
$ws.Range("K234").Value = 500
$ws.Range("L234").Value = 500
$ws.Range("M234").Value = 500
$ws.Range("P234").Value = 500

$ws.Range("D235").Value = 44285
$ws.Range("J235").Value = 3000
$ws.Range("K235").Value = 800
$ws.Range("L235").Value = 800
$ws.Range("M235").Value = 800
$ws.Range("O235").Value = "Región del Maule"
$ws.Range("P235").Value = 800

$ws.Range("D236").Value = 44498
$ws.Range("J236").Value = 6000
$ws.Range("K236").Value = 600
$ws.Range("L236").Value = 600
$ws.Range("M236").Value = 600
$ws.Range("O236").Value = "Provincia del Elquí"
$ws.Range("P236").Value = 600

$ws.Range("D237").Value = 44179
$ws.Range("J237").Value = 3000
$ws.Range("K237").Value = 700
$ws.Range("L237").Value = 700
$ws.Range("M237").Value = 700
$ws.Range("P237").Value = 700

$ws.Range("I238").Value = "Primera"
$ws.Range("J238").Value = 2000
$ws.Range("K238").Value = 500
$ws.Range("L238").Value = 500
$ws.Range("M238").Value = 500
$ws.Range("P238").Value = 500

$ws.Range("D239").Value = 44418
$ws.Range("J239").Value = 3000
$ws.Range("K239").Value = 350
$ws.Range("L239").Value = 350
$ws.Range("M239").Value = 350
$ws.Range("P239").Value = 350

$ws.Range("D240").Value = 44595
$ws.Range("K240").Value = 800
$ws.Range("L240").Value = 800
$ws.Range("M240").Value = 800
$ws.Range("P240").Value = 800

$ws.Range("D241").Value = 44628
$ws.Range("I241").Value = "Segunda"
$ws.Range("J241").Value = 2000
$ws.Range("K241").Value = 1000
$ws.Range("L241").Value = 1000
$ws.Range("M241").Value = 1000
$ws.Range("P241").Value = 1000

$ws.Range("D242").Value = 44335
$ws.Range("K242").Value = 550
$ws.Range("L242").Value = 550
$ws.Range("M242").Value = 550
$ws.Range("P242").Value = 550

$ws.Range("D243").Value = 44552
$ws.Range("J243").Value = 3000
$ws.Range("K243").Value = 600
$ws.Range("L243").Value = 600
$ws.Range("M243").Value = 600
$ws.Range("P243").Value = 600

$ws.Range("D244").Value = 44544
$ws.Range("J244").Value = 4000
$ws.Range("K244").Value = 500
$ws.Range("L244").Value = 500
$ws.Range("M244").Value = 500
$ws.Range("P244").Value = 500

$ws.Range("A245").Value = 5
$ws.Range("B245").Value = "Macroferia Regional de Talca"
$ws.Range("C245").Value = "Maule"
$ws.Range("D245").Value = 44160
$ws.Range("D245").NumberFormat = $ws.Range("D244").NumberFormat
$ws.Range("E245").Value = 7
$ws.Range("F245").Value = 100112006
$ws.Range("G245").Value = "Repollo"
$ws.Range("H245").Value = "Crespo record"
$ws.Range("I245").Value = "Primera"
$ws.Range("J245").Value = 2000
$ws.Range("K245").Value = 900
$ws.Range("L245").Value = 900
$ws.Range("M245").Value = 900
$ws.Range("N245").Value = "$/unidad"
$ws.Range("O245").Value = "Región del Maule"
$ws.Range("P245").Value = 900
$ws.Range("Q245").Value = 1
$ws.Range("R245").Value = "Hortaliza"
